# Edit: Remove renewables from BAU guaranteed dispatch and set coal to bid
# at its expected capacity factor.

$wb = $excel.ActiveWorkbook

# --- BDSBaPCF sheet: set "hard coal" to bid at Expected Capacity Factor (0) ---
$wsData = $wb.Worksheets.Item("BDSBaPCF")
$wsData.Range("B2").Value = 0

# --- About sheet: append explanatory notes about coal ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A24").Value = "For the United States, we have set coal to 0 as of version 3.4. This reflects"
$wsAbout.Range("A25").Value = "the fact that certain air quality / environmental restrictions, as well as current"
$wsAbout.Range("A26").Value = "supply chain logistics, limit the amount the coal dispatches annually. "

$wb.Save()
